# "updated forms with labs"
#
# This script updates the rarrive.xlsx XLSForm workbook:
#   - survey sheet: renames the "region" select_one label/hint from
#     "Region" / quick search('regions') to "Facility Lab" / quick search('labs')
#   - choices sheet: renames the region list's key/label columns from
#     region_key/region to lab_key/lab, and removes the obsolete "rider"
#     choice-list row (shifting the rows below it up by one)
#   - settings sheet: adds a "version" column with a form version number

$wb = $excel.ActiveWorkbook

$survey   = $wb.Worksheets.Item("survey")
$choices  = $wb.Worksheets.Item("choices")
$settings = $wb.Worksheets.Item("settings")

# ---------------------------------------------------------------------------
# survey sheet: "select_one region" row -> relabel as the facility lab prompt
# ---------------------------------------------------------------------------
$survey.Range("C6").Value = "Facility Lab"
$survey.Range("E6").Value = "quick search('labs')"

# ---------------------------------------------------------------------------
# choices sheet: rename the region list's key/label fields to lab_key/lab,
# then drop the "rider" choice-list row, shifting rows below it up by one.
# ---------------------------------------------------------------------------
$choices.Range("B2").Value = "lab_key"
$choices.Range("C2").Value = "lab"

# Remove the "rider" row (row 4): shift the following rows (stype,
# condition) up by one, and clear out what is now the trailing row.
$choices.Range("A4").Value = $choices.Range("A5").Value2
$choices.Range("B4").Value = $choices.Range("B5").Value2
$choices.Range("C4").Value = $choices.Range("C5").Value2

$choices.Range("A5").Value = $choices.Range("A6").Value2
$choices.Range("B5").Value = $choices.Range("B6").Value2
$choices.Range("C5").Value = $choices.Range("C6").Value2

$choices.Range("A6").ClearContents()
$choices.Range("B6").ClearContents()
$choices.Range("C6").ClearContents()

# ---------------------------------------------------------------------------
# settings sheet: add a "version" column with the form's version number
# ---------------------------------------------------------------------------
$settings.Range("C1").Value = "version"
$settings.Range("C1").Font.Bold = $true

$settings.Range("C2").Value = 2016041301
$settings.Range("C2").HorizontalAlignment = -4152
